$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark (currently a zero-width
#        bookmark sitting right after "All staff are permitted to do image
#        uploads.  " at the end of business rule 017's body paragraph). ---
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- 2. Rule "017 - Upload Authorization": the run holding the two literal
#        spaces right after the bold title (before the run of tab stops
#        that lead up to "Date: 10/02/2019") is deleted outright. ---
$r017 = $d.Content
$r017.Find.Execute("Upload Authorization", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$gap017 = $d.Range($r017.End, $r017.End + 2)
$gap017.Delete()

# --- 3. Rule "015 - Showing Appointment Travel Time": same two-space run
#        right after the bold title is deleted. ---
$r015 = $d.Content
$r015.Find.Execute("Travel Time", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$gap015 = $d.Range($r015.End, $r015.End + 2)
$gap015.Delete()

# --- 4. The "_GoBack" bookmark re-appears inside rule 015's header line,
#        now sitting between the first and second tab stops that follow
#        the title (title -> tab -> [_GoBack] -> tab -> tab -> tab -> tab
#        -> "Date: 10/02/2019"). ---
$r015b = $d.Content
$r015b.Find.Execute("Travel Time", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkPos = $r015b.End + 1
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
